$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: quarter/period headers (D8:M8)
$ws.Range("D8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("J8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# Row 9: publish-date headers (D9:M9)
$ws.Range("D9").Value = "1400-10-29 (3)"
$ws.Range("E9").Value = "1401-04-12 (11)"
$ws.Range("F9").Value = "1401-04-28 (4)"
$ws.Range("G9").Value = "1401-09-14 (4)"
$ws.Range("H9").Value = "1401-10-28 (2)"
$ws.Range("I9").Value = "1402-02-25 (8)"
$ws.Range("J9").Value = "1401-04-28"
$ws.Range("K9").Value = "1401-09-14 (2)"
$ws.Range("L9").Value = "1401-10-28"
$ws.Range("M9").Value = "1402-02-25"

# Data rows
# Row 11
$ws.Range("D11").Value = 11525
$ws.Range("E11").Value = 17995
$ws.Range("F11").Value = 14529
$ws.Range("G11").Value = 16400
$ws.Range("H11").Value = 15290
$ws.Range("I11").Value = 18831
$ws.Range("J11").Value = 19124
$ws.Range("K11").Value = 21421
$ws.Range("L11").Value = 21070
$ws.Range("M11").Value = 16889

# Row 12
$ws.Range("D12").Value = -4735
$ws.Range("E12").Value = -11137
$ws.Range("F12").Value = -7222
$ws.Range("G12").Value = -8142
$ws.Range("H12").Value = -8791
$ws.Range("I12").Value = -10513
$ws.Range("J12").Value = -10473
$ws.Range("K12").Value = -11005
$ws.Range("L12").Value = -12038
$ws.Range("M12").Value = -9880

# Row 13
$ws.Range("D13").Value = 6789
$ws.Range("E13").Value = 6858
$ws.Range("F13").Value = 7308
$ws.Range("G13").Value = 8257
$ws.Range("H13").Value = 6499
$ws.Range("I13").Value = 8318
$ws.Range("J13").Value = 8651
$ws.Range("K13").Value = 10416
$ws.Range("L13").Value = 9033
$ws.Range("M13").Value = 7009

# Row 14
$ws.Range("D14").Value = -1595
$ws.Range("E14").Value = -2519
$ws.Range("F14").Value = -2819
$ws.Range("G14").Value = -2905
$ws.Range("H14").Value = -2153
$ws.Range("I14").Value = -3833
$ws.Range("J14").Value = -3203
$ws.Range("K14").Value = -3768
$ws.Range("L14").Value = -3849
$ws.Range("M14").Value = -2901

# Row 17
$ws.Range("D17").Value = 5194
$ws.Range("E17").Value = 4339
$ws.Range("F17").Value = 4489
$ws.Range("G17").Value = 5352
$ws.Range("H17").Value = 4346
$ws.Range("I17").Value = 4484
$ws.Range("J17").Value = 5448
$ws.Range("K17").Value = 6648
$ws.Range("L17").Value = 5184
$ws.Range("M17").Value = 4108

# Row 18
$ws.Range("D18").Value = -479
$ws.Range("E18").Value = -603
$ws.Range("F18").Value = -637
$ws.Range("G18").Value = -609
$ws.Range("H18").Value = -604
$ws.Range("I18").Value = -707
$ws.Range("J18").Value = -845
$ws.Range("K18").Value = -797
$ws.Range("L18").Value = -769
$ws.Range("M18").Value = -607

# Row 19
$ws.Range("D19").Value = 158
$ws.Range("E19").Value = -32
$ws.Range("F19").Value = -15
$ws.Range("G19").Value = 100
$ws.Range("H19").Value = 35
$ws.Range("I19").Value = 236
$ws.Range("J19").Value = -167
$ws.Range("K19").Value = -22
$ws.Range("L19").Value = -322
$ws.Range("M19").Value = 592

# Row 20
$ws.Range("D20").Value = 4873
$ws.Range("E20").Value = 3705
$ws.Range("F20").Value = 3837
$ws.Range("G20").Value = 4843
$ws.Range("H20").Value = 3778
$ws.Range("I20").Value = 4014
$ws.Range("J20").Value = 4436
$ws.Range("K20").Value = 5830
$ws.Range("L20").Value = 4093
$ws.Range("M20").Value = 4093

# Row 21
$ws.Range("D21").Value = -406
$ws.Range("E21").Value = 136
$ws.Range("F21").Value = -142
$ws.Range("G21").Value = -354
$ws.Range("H21").Value = -310
$ws.Range("I21").Value = -15
$ws.Range("J21").Value = -311
$ws.Range("K21").Value = -319
$ws.Range("L21").Value = -250
$ws.Range("M21").Value = 47

# Row 22
$ws.Range("D22").Value = 4467
$ws.Range("E22").Value = 3840
$ws.Range("F22").Value = 3694
$ws.Range("G22").Value = 4489
$ws.Range("H22").Value = 3468
$ws.Range("I22").Value = 3999
$ws.Range("J22").Value = 4126
$ws.Range("K22").Value = 5511
$ws.Range("L22").Value = 3843
$ws.Range("M22").Value = 4140

# Row 24
$ws.Range("D24").Value = 4467
$ws.Range("E24").Value = 3840
$ws.Range("F24").Value = 3694
$ws.Range("G24").Value = 4489
$ws.Range("H24").Value = 3468
$ws.Range("I24").Value = 3999
$ws.Range("J24").Value = 4126
$ws.Range("K24").Value = 5511
$ws.Range("L24").Value = 3843
$ws.Range("M24").Value = 4140

# Row 26
$ws.Range("D26").Value = 5520
$ws.Range("E26").Value = 6141
$ws.Range("F26").Value = 6430
$ws.Range("G26").Value = 5747
$ws.Range("H26").Value = 5268
$ws.Range("I26").Value = 5446
$ws.Range("J26").Value = 5105
$ws.Range("K26").Value = 12897
$ws.Range("L26").Value = 11513
$ws.Range("M26").Value = 8800
